# Add 2022-Q3 data:
#  1) Insert a new worksheet "2022-Q3" right after "总计" (and right before "2022-Q2"),
#     populated with the fund-holding detail rows for that quarter.
#  2) Update the "总计" (summary) sheet: insert a new row right below the header
#     with the 2022-Q3 totals, shifting every existing data row down by one and
#     renumbering the running index column (A).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2    = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q3" detail sheet, placed right before "2022-Q2"
# ---------------------------------------------------------------------------
# Duplicate the "2022-Q2" sheet (rather than Worksheets.Add a blank one) so
# the new sheet inherits all of its formatting (bold/centered header style,
# borders, column widths, etc.) verbatim. The duplicate is placed right
# before "2022-Q2", exactly where "2022-Q3" belongs.
$wsQ2.Copy($wsQ2, $null)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# The source sheet has 10 data rows (rows 2-11); 2022-Q3 only needs 6 data
# rows (rows 2-7), so drop the extra ones (this also shrinks the sheet's
# dimension/used-range back down to A1:H7).
$wsQ3.Range("A8:H11").Delete()

# Header values (overwrite the copied "2022-Q2" header text with itself /
# the same captions - kept explicit for clarity and resilience).
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# Columns that must stay plain numbers
$numCols = @("A", "H")
# Columns that must stay text even though they look numeric (fund codes,
# percentages stored as strings, etc.)
$textCols = @("B", "C", "D", "E", "F")

foreach ($col in $textCols) {
    $wsQ3.Range($col + "1:" + $col + "7").NumberFormat = "@"
}

function Set-Q3Row($r, $idx, $code, $fundName, $scale, $stockPos, $posPct, $heldValue, $isHeldValueText, $rank) {
    $wsQ3.Range("A$r").Value = $idx
    $wsQ3.Range("B$r").Value = $code
    $wsQ3.Range("C$r").Value = $fundName
    $wsQ3.Range("D$r").Value = $scale
    $wsQ3.Range("E$r").Value = $stockPos
    $wsQ3.Range("F$r").Value = $posPct
    if ($isHeldValueText) {
        $wsQ3.Range("G$r").NumberFormat = "@"
        $wsQ3.Range("G$r").Value = $heldValue
    } else {
        $wsQ3.Range("G$r").NumberFormat = "General"
        $wsQ3.Range("G$r").Value = $heldValue
    }
    $wsQ3.Range("H$r").Value = $rank
}

Set-Q3Row 2 0 "006323" "合煦智远嘉选混合A"      "0.78" "73.07" "6.99" "0.0545" $true  1
Set-Q3Row 3 1 "673090" "西部利得个股精选股票A"   "1.25" "86.69" "2.47" "0.0309" $true  6
Set-Q3Row 4 2 "013262" "西部利得个股精选股票C"   "0.95" "86.69" "2.47" "0.0235" $true  6
Set-Q3Row 5 3 "013204" "恒生前海恒源天利债A"     "1.30" "32.00" "1.15" "0.0150" $true  10
Set-Q3Row 6 4 "006324" "合煦智远嘉选混合C"       "0.14" "73.07" "6.99" "0.0098" $true  1
Set-Q3Row 7 5 "013205" "恒生前海恒源天利债C"     "0.00" "32.00" "1.15" 0        $false 10

$wsQ3.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet
# ---------------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

# Bring formatting for the new row in from the row right below it (which used
# to be row 2, now shifted to row 3) so styles stay identical to the rest of
# the table.
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 6
$wsTotal.Range("D2").Value = 0.13

# Renumber the running index (column A) for all the rows that shifted down.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5
$wsTotal.Range("A8").Value = 6
$wsTotal.Range("A9").Value = 7

$wsTotal.Range("A1").Select()

$wb.Save()
